# SkillRef.xlsx edit
#
# Per the commit ("unify the conception of DataNode, DataTable, Entity")
# the only real content-level edit in this workbook is the sheet being
# renamed from "Property1" to "DataNode", and the surviving cursor
# selection moving to C38 (the cell the author was on when the sheet was
# re-saved). Everything else in the source diff (fileVersion/rupBuild,
# bookViews, xr:*/x15 namespace & uid churn, the absPath machine path,
# default row height / baseColWidth / sub-pixel column-width shifts, the
# phoneticPr + extra 9pt font, the "Normal"->"常规" cell-style display
# name, the x15:timelineStyles ext) is Excel/OS/locale boilerplate that
# gets stamped automatically whenever the file is opened & resaved by a
# different Excel build (here: Mac Excel -> Windows Excel, en-US -> zh-CN)
# and is not something reachable - or meaningful to force - through the
# COM object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Property1" -> "DataNode"
$ws.Name = "DataNode"

# Selection left on C38 (was A9).
[void]$ws.Range("C38").Select()
